# Apply VerveStacks CHE grids update (2025-08-09 16:56) to the "existing_stock" sheet.
# All real content changes in this revision are confined to sheet "existing_stock";
# "weo_pg" and "ccs_retrofits" only had internal revision-tracking ids touched, which
# carry no user-visible data and are regenerated automatically by Excel on save.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("existing_stock")

$ws.Range("E14").Value2 = 0.086
$ws.Range("H14").Value2 = 60.50000000000001
$ws.Range("E15").Value2 = 0.066
$ws.Range("H15").Value2 = 60.500000000000014
$ws.Range("C82").Value2 = 'e_w391576135-220'
$ws.Range("E82").Value2 = 0.061
$ws.Range("G82").Value2 = 2783.0
$ws.Range("H82").Value2 = 66.55000000000001
$ws.Range("I82").Value2 = 2.8875
$ws.Range("C83").Value2 = 'e_w1284913429-220'
$ws.Range("E83").Value2 = 0.05
$ws.Range("G83").Value2 = 3267.0000000000005
$ws.Range("H83").Value2 = 78.65
$ws.Range("I83").Value2 = 3.1500000000000004
$ws.Range("E96").Value2 = 0.0012
$ws.Range("C97").Value2 = 'e_w27435934-220'
$ws.Range("E97").Value2 = 0.0017
$ws.Range("G97").Value2 = 1336.4999999999998
$ws.Range("C98").Value2 = 'e_w234983117-220'
$ws.Range("E99").Value2 = 0.001
$ws.Range("E100").Value2 = 0.0011
$ws.Range("G100").Value2 = 1336.5
$ws.Range("C101").Value2 = 'e_w83861269-220'
$ws.Range("E101").Value2 = 0.0018
$ws.Range("G101").Value2 = 1336.5000000000002
$ws.Range("H101").Value2 = 21.450000000000003
$ws.Range("C102").Value2 = 'e_CH17-380'
$ws.Range("E102").Value2 = 0.0012
$ws.Range("C103").Value2 = 'e_w281809991-220'
$ws.Range("E103").Value2 = 0.0012
$ws.Range("C104").Value2 = 'e_w97941869-220'
$ws.Range("E104").Value2 = 0.0015
$ws.Range("C105").Value2 = 'e_CH60-225'
$ws.Range("E105").Value2 = 0.0034
$ws.Range("G105").Value2 = 1336.5000000000002
$ws.Range("C106").Value2 = 'e_w33271433-220'
$ws.Range("E106").Value2 = 0.003
$ws.Range("C107").Value2 = 'e_w127004407-380'
$ws.Range("E107").Value2 = 0.0017
$ws.Range("G107").Value2 = 1336.4999999999998
$ws.Range("C108").Value2 = 'e_CH31-220'
$ws.Range("E108").Value2 = 0.001
$ws.Range("G108").Value2 = 1336.5
$ws.Range("C109").Value2 = 'e_w234983117-220'
$ws.Range("E109").Value2 = 0.0023
$ws.Range("G109").Value2 = 1336.5
$ws.Range("H109").Value2 = 21.450000000000006
$ws.Range("C110").Value2 = 'e_w89977424-220'
$ws.Range("E110").Value2 = 0.0015
$ws.Range("C111").Value2 = 'e_w1105061707-220'
$ws.Range("E111").Value2 = 0.0016
$ws.Range("C112").Value2 = 'e_CH57-220'
$ws.Range("E112").Value2 = 0.0058000000000000005
$ws.Range("C113").Value2 = 'e_w190819048-220'
$ws.Range("E113").Value2 = 0.001
$ws.Range("C114").Value2 = 'e_w281809991-220'
$ws.Range("E114").Value2 = 0.0013
$ws.Range("C115").Value2 = 'e_CH31-220'
$ws.Range("E115").Value2 = 0.0012
$ws.Range("C116").Value2 = 'e_CH60-225'
$ws.Range("E116").Value2 = 0.0011
$ws.Range("C117").Value2 = 'e_w109037817-220'
$ws.Range("E117").Value2 = 0.0021000000000000003
$ws.Range("G117").Value2 = 1336.5
$ws.Range("C118").Value2 = 'e_CH11-220'
$ws.Range("E118").Value2 = 0.0070999999999999995
$ws.Range("G118").Value2 = 1336.5000000000002
$ws.Range("C119").Value2 = 'e_w1105061707-220'
$ws.Range("E119").Value2 = 0.001
$ws.Range("C121").Value2 = 'e_w281809991-220'
$ws.Range("E121").Value2 = 0.005
$ws.Range("C122").Value2 = 'e_w97941869-220'
$ws.Range("E122").Value2 = 0.0015
$ws.Range("C123").Value2 = 'e_r5378910-220'
$ws.Range("E123").Value2 = 0.0013
$ws.Range("C124").Value2 = 'e_w89977424-220'
$ws.Range("E124").Value2 = 0.0011
$ws.Range("C125").Value2 = 'e_CH11-220'
$ws.Range("E125").Value2 = 0.0012
$ws.Range("C126").Value2 = 'e_CH17-380'
$ws.Range("E126").Value2 = 0.0021000000000000003
$ws.Range("C127").Value2 = 'e_w165254212-220'
$ws.Range("E127").Value2 = 0.0058
$ws.Range("C128").Value2 = 'e_w1105061707-220'
$ws.Range("E128").Value2 = 0.0045
$ws.Range("C130").Value2 = 'e_w97941869-220'
$ws.Range("E130").Value2 = 0.0022
$ws.Range("C131").Value2 = 'e_r5378910-220'
$ws.Range("E131").Value2 = 0.0014
$ws.Range("C132").Value2 = 'e_CH31-220'
$ws.Range("E132").Value2 = 0.0108
$ws.Range("H132").Value2 = 21.45
$ws.Range("C133").Value2 = 'e_w127004407-380'
$ws.Range("E133").Value2 = 0.004900000000000001
$ws.Range("H133").Value2 = 21.450000000000003
$ws.Range("E135").Value2 = 0.21731537653220406
$ws.Range("E136").Value2 = 0.1663793251004252
$ws.Range("E137").Value2 = 0.19782569372870323
$ws.Range("E138").Value2 = 0.21077006448261207
$ws.Range("E139").Value2 = 0.16277902359433066
$ws.Range("E140").Value2 = 0.16152827258311295
$ws.Range("E141").Value2 = 0.1371001651339535
$ws.Range("E142").Value2 = 0.1586174511333161
$ws.Range("E143").Value2 = 0.15521278721895346
$ws.Range("E144").Value2 = 0.16554834018408843
$ws.Range("E145").Value2 = 0.15291072157643879
$ws.Range("E146").Value2 = 0.1727926292604506
$ws.Range("E147").Value2 = 0.1930981018275324
$ws.Range("E148").Value2 = 0.21471510601685545
$ws.Range("E149").Value2 = 0.18308354646436523
$ws.Range("E150").Value2 = 0.19614947844032105
$ws.Range("E151").Value2 = 0.1533817713118708
$ws.Range("E152").Value2 = 0.19328994063107527
$ws.Range("E153").Value2 = 0.1982862967966156
$ws.Range("E154").Value2 = 0.15409499379434963
$ws.Range("E155").Value2 = 0.16699472878703805
$ws.Range("E156").Value2 = 0.13606784444360143
$ws.Range("E158").Value2 = 0.20091315882928704
$ws.Range("E159").Value2 = 0.211523174241075
